$wb = $excel.ActiveWorkbook

# Sheet1: new leaderboard entry "j" with score 0
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A15").Value = "j"
$ws1.Range("B15").Value = 0

# Sheet4: new leaderboard entry "j" with score 2084
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("A4").Value = "j"
$ws4.Range("B4").Value = 2084

# Sheet6: new leaderboard entry "l7" with score 0
$ws6 = $wb.Worksheets.Item("Sheet6")
$ws6.Range("A4").Value = "l7"
$ws6.Range("B4").Value = 0
